$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.246.51'
$ws.Range("E2").Value = '  +0.20%  '

$ws.Range("D3").Value = '1.862.47'
$ws.Range("E3").Value = '  +0.64%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7103'
$ws.Range("E5").Value = '  +0.49%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '238.09'
$ws.Range("E6").Value = '  -0.42%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.08155'
$ws.Range("E8").Value = '  +9.76%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3044'
$ws.Range("E9").Value = '  -0.53%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.32'
$ws.Range("E10").Value = '  -1.17%  '

$ws.Range("E11").Value = '  +0.45%  '

$ws.Range("D12").Value = '1.865.82'
$ws.Range("E12").Value = '  +0.86%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.183'
$ws.Range("E13").Value = '  -1.01%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7108'
$ws.Range("E14").Value = '  -2.47%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.84'
$ws.Range("E15").Value = '  +0.97%  '

$ws.Range("D16").Value = '29.257.62'
$ws.Range("E16").Value = '  +0.41%  '

$ws.Range("B17").Value = 'Uniswap'
$ws.Range("C17").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.794'
$ws.Range("E17").Value = '  +0.34%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007892'
$ws.Range("E18").Value = '  +2.99%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.42'
$ws.Range("E19").Value = '  +2.43%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '237.45'
$ws.Range("E20").Value = '  -0.67%  '

$ws.Range("E21").Value = '  +0.10%  '

$ws.Range("D22").Value = '2.101.69'
$ws.Range("E22").Value = '  +1.32%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  +0.04%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.422'
$ws.Range("E24").Value = '  -2.52%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '162.68'
$ws.Range("E25").Value = '  +0.96%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.975'
$ws.Range("E26").Value = '  -0.51%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1457'
$ws.Range("E27").Value = '  -0.16%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.12'
$ws.Range("E28").Value = '  +0.03%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.965'
$ws.Range("E29").Value = '  -0.69%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.427'
$ws.Range("E30").Value = '  +2.07%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.486'
$ws.Range("E31").Value = '  -0.35%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.400'
$ws.Range("E32").Value = '  -3.20%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.039'
$ws.Range("E33").Value = '  +1.14%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05226'
$ws.Range("E34").Value = '  +0.45%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.172'
$ws.Range("E35").Value = '  -1.44%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7091'
$ws.Range("E36").Value = '  +0.59%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9996'
$ws.Range("E37").Value = '  -3.13%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.671'
$ws.Range("E38").Value = '  +0.40%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01857'
$ws.Range("E39").Value = '  -0.83%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.733'
$ws.Range("E40").Value = '  +2.09%  '

$ws.Range("D41").Value = '1.146.68'
$ws.Range("E41").Value = '  +6.70%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9243'
$ws.Range("E42").Value = '  -2.51%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4289'
$ws.Range("E43").Value = '  -0.22%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.886'
$ws.Range("E44").Value = '  -2.23%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '70.39'
$ws.Range("E45").Value = '  -0.02%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9997'
$ws.Range("E46").Value = '  -0.02%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.86'
$ws.Range("E47").Value = '  -0.02%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.777'
$ws.Range("E48").Value = '  +1.77%  '

$ws.Range("D49").Value = '2.001.15'
$ws.Range("E49").Value = '  +1.04%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.236'
$ws.Range("E50").Value = '  +1.39%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.969'
$ws.Range("E51").Value = '  -1.57%  '
